$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: A2 FAPs(unchanged), B2 -> Nlgn1, C2 -> Nrxn1, D2 FAPs(unchanged)
$ws.Cells.Item(2,2).Value = "Nlgn1"
$ws.Cells.Item(2,3).Value = "Nrxn1"
$ws.Cells.Item(2,9).Value = 1
$ws.Cells.Item(2,10).Value = 1
$ws.Cells.Item(2,15).Value = 0.8459226744507667
$ws.Cells.Item(2,16).Value = 0.8459226744507669
$ws.Cells.Item(2,19).Value = 0.8459226744507667
$ws.Cells.Item(2,20).Value = 0.8459226744507669

# Row 3: A3 FAPs(unchanged), B3 -> Nlgn1, C3 -> Nrxn1, D3 -> MuSCs
$ws.Cells.Item(3,2).Value = "Nlgn1"
$ws.Cells.Item(3,3).Value = "Nrxn1"
$ws.Cells.Item(3,4).Value = "MuSCs"
$ws.Cells.Item(3,9).Value = 1
$ws.Cells.Item(3,10).Value = 1
$ws.Cells.Item(3,13).Value = 0.03648100000000001
$ws.Cells.Item(3,14).Value = 0.109443
$ws.Cells.Item(3,15).Value = 0.1407969268413801
$ws.Cells.Item(3,16).Value = 0.1407969268413801
$ws.Cells.Item(3,17).Value = 0.0005504982900000001
$ws.Cells.Item(3,18).Value = 0.00495448461
$ws.Cells.Item(3,19).Value = 0.1407969268413801
$ws.Cells.Item(3,20).Value = 0.1407969268413801

# Row 4: was "MuSCs -> FAPs"; becomes a new pair "FAPs -> Resolving-Mac" with new TPM-derived numbers
$ws.Cells.Item(4,1).Value = "FAPs"
$ws.Cells.Item(4,2).Value = "Nlgn1"
$ws.Cells.Item(4,3).Value = "Nrxn1"
$ws.Cells.Item(4,4).Value = "Resolving-Mac"
$ws.Cells.Item(4,7).Value = 0.01509
$ws.Cells.Item(4,8).Value = 0.04527
$ws.Cells.Item(4,9).Value = 1
$ws.Cells.Item(4,10).Value = 1
$ws.Cells.Item(4,11).Value = 2
$ws.Cells.Item(4,12).Value = 0.6666666666666666
$ws.Cells.Item(4,13).Value = 0.003441
$ws.Cells.Item(4,14).Value = 0.010323
$ws.Cells.Item(4,15).Value = 0.0132803987078531
$ws.Cells.Item(4,16).Value = 0.0132803987078531
$ws.Cells.Item(4,17).Value = 0.00005192469
$ws.Cells.Item(4,18).Value = 0.00046732221
$ws.Cells.Item(4,19).Value = 0.0132803987078531
$ws.Cells.Item(4,20).Value = 0.0132803987078531

# Row 5 removed entirely
$ws.Rows.Item(5).Delete()
